$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F17").Value = '18_hazards_to_humans_and_domestic_animals'
$ws.Range("F19").Value = 'ppe'
$ws.Range("F20").Value = 'ppe'
$ws.Range("F21").Value = 'ppe'
$ws.Range("F26").Value = 'application instructions || env warning - species'
$ws.Range("F27").Value = 'env warning - water'
$ws.Range("F29").Value = 'env warning - water || off target movement'
$ws.Range("F32").Value = '32_physical_and_chemical_hazards'
$ws.Range("F38").Value = 'application instructions'
$ws.Range("F39").Value = 'application instructions'
$ws.Range("F40").Value = 'application instructions'
$ws.Range("F41").Value = 'application instructions'
$ws.Range("F42").Value = 'application instructions'
$ws.Range("F46").Value = 'use restrictions'
$ws.Range("F47").Value = 'use restrictions'
$ws.Range("F48").Value = 'use restrictions'
$ws.Range("F49").Value = 'use restrictions'
$ws.Range("F51").Value = '135_product_information'
$ws.Range("F52").Value = '135_product_information'
$ws.Range("F53").Value = '135_product_information'
$ws.Range("F61").Value = 'application instructions'
$ws.Range("F62").Value = 'application instructions'
$ws.Range("F64").Value = 'application instructions'
$ws.Range("F65").Value = 'application instructions'
$ws.Range("F66").Value = 'application instructions'
$ws.Range("F67").Value = 'application instructions'
$ws.Range("F69").Value = 'application instructions'
$ws.Range("F70").Value = 'mixing'
$ws.Range("F72").Value = 'mixing'
$ws.Range("F74").Value = 'chemigation'
$ws.Range("F75").Value = 'irrigation || application instructions || chemigation'
$ws.Range("F76").Value = 'irrigation'
$ws.Range("F77").Value = 'safety procedures'
$ws.Range("F78").Value = 'chemigation'
$ws.Range("F79").Value = 'chemigation'
$ws.Range("F85").Value = 'off target movement'
$ws.Range("F86").Value = 'off target movement'
$ws.Range("F87").Value = 'off target movement'
$ws.Range("F88").Value = 'off target movement'
$ws.Range("F90").Value = 'off target movement'
$ws.Range("F92").Value = 'off target movement'
$ws.Range("F93").Value = 'off target movement'
$ws.Range("F94").Value = 'off target movement'
$ws.Range("F95").Value = 'off target movement'
$ws.Range("F98").Value = 'off target movement'
$ws.Range("F100").Value = 'off target movement'
$ws.Range("F109").Value = 'off target movement'
$ws.Range("F114").Value = 'off target movement'
$ws.Range("F115").Value = 'off target movement'
$ws.Range("F117").Value = 'application instructions'
$ws.Range("F120").Value = 'application instructions'
$ws.Range("F122").Value = 'application instructions'
$ws.Range("F123").Value = 'application instructions'
$ws.Range("F125").Value = 'application instructions'
$ws.Range("F126").Value = 'application instructions'
$ws.Range("F128").Value = 'mixing'
$ws.Range("F130").Value = 'mixing'
$ws.Range("F131").Value = 'mixing'
$ws.Range("F196").Value = 'mixing'
$ws.Range("F250").Value = 'use restrictions'
$ws.Range("F253").Value = '154_pesticide_storage'
$ws.Range("F255").Value = '154_pesticide_storage'
